$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "2021年"
$ws.Range("A11").Style = $ws.Range("A10").Style
$ws.Range("B11").Value = 123.04
$ws.Range("C11").Value = 20.62
$ws.Range("D11").Value = 3.95
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = 127.18
$ws.Range("G11").Value = 558.33
$ws.Range("H11").Value = 95.68000000000001
$ws.Range("I11").Value = 111.5
$ws.Range("J11").Value = 16.17
$ws.Range("K11").Value = 6476.34
$ws.Range("L11").Value = 22.97
$ws.Range("M11").Value = 17.56
$ws.Range("N11").Value = 3.07
$ws.Range("O11").Value = 37.49
$ws.Range("P11").Value = 458.69
$ws.Range("Q11").Value = 40.99
$ws.Range("R11").Value = 15.47
$ws.Range("S11").Value = 82.83
$ws.Range("T11").Value = 89.04000000000001
$ws.Range("U11").Value = -41.22
$ws.Range("V11").Value = -83.81
$ws.Range("W11").Value = 791.01
$ws.Range("X11").Value = 40.42
$ws.Range("Y11").Value = 1462.1
$ws.Range("Z11").Value = 245.21
$ws.Range("AA11").Value = 24.05
$ws.Range("AB11").Value = 343.16
$ws.Range("AC11").Value = 154.43
$ws.Range("AD11").Value = 128.16
$ws.Range("AE11").Value = 45.49
$ws.Range("AF11").Value = 311.58
$ws.Range("AG11").Value = 88.61
$ws.Range("AH11").Value = 111.61
$ws.Range("AI11").Value = -3.5
$ws.Range("AJ11").Value = 8.6
$ws.Range("AK11").Value = 127.37
$ws.Range("AL11").Value = 29.73
$ws.Range("AM11").Value = 184.13
$ws.Range("AN11").Value = 11.72
$ws.Range("AO11").Value = 33.74
$ws.Range("AP11").Value = 491.21
$ws.Range("AQ11").Value = 147.88
